# Append two new trading-log rows (54 and 55) to the active sheet,
# mirroring the TRADING_ATTEMPT / POSITION_OPENED pair that the bot
# logs for each trade.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 54: TRADING_ATTEMPT
$ws.Range("A54").Value = "2025-09-29T01:28:00.101981"
$ws.Range("B54").Value = "TRADING_ATTEMPT"
$ws.Range("C54").Value = "BTC"
$ws.Range("D54").Value = "UNKNOWN"
$ws.Range("E54").Value = 112265.3486426338
$ws.Range("K54").Value = "ATTEMPT"
$ws.Range("L54").Value = "Attempting trade 1/1"

# Row 55: POSITION_OPENED
$ws.Range("A55").Value = "2025-09-29T01:28:01.747841"
$ws.Range("B55").Value = "POSITION_OPENED"
$ws.Range("C55").Value = "BTC"
$ws.Range("D55").Value = "UNKNOWN"
$ws.Range("E55").Value = 112265.3486426338
$ws.Range("F55").Value = 3600
$ws.Range("G55").Value = 40
$ws.Range("H55").Value = 0.2775619584775681
$ws.Range("K55").Value = "SUCCESS"
